# fix: implement export menu Fixes #7
#
# Adds a new membership record ("Kevin Smith" / 123 Adventure Drive) as a
# new row 2 on the "LCRA membership" sheet (pushing the existing record
# down to row 3), and turns the new record's e-mail address into a
# mailto: hyperlink.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a fresh blank row above the current data row (old row 2 becomes
# row 3, carrying all of its values/formatting down with it).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Select() | Out-Null

# --- New member's address -------------------------------------------------
$ws.Range("A2").Value = "123 Adventure Drive"
$ws.Range("B2").Value = 123
$ws.Range("C2").Value = "Adventure Drive"
$ws.Range("D2").Value = "A0A0A0"

# --- New member's primary contact info -------------------------------------
$ws.Range("E2").Value = "Kevin"
$ws.Range("F2").Value = "Smith"
$ws.Range("G2").Value = "ksmith@email.com"
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

# --- Record metadata ---------------------------------------------------
$ws.Range("AK2").Value = "2023-02-23T03:19:09+00:00"
$ws.Range("AL2").Value = "testuser"
$ws.Range("AM2").Value = "2nd entry"

# --- Membership-year history ---------------------------------------------
$ws.Range("AN2").Value = 0
$ws.Range("AO2").Value = ""
$ws.Range("AP2").Value = ""
$ws.Range("AQ2").Value = ""
$ws.Range("AR2").Value = 0
$ws.Range("AS2").Value = 1
$ws.Range("AT2").Value = "Corn Roast"
$ws.Range("AU2").Value = "2023-06-11T16:00:00-04:00"
$ws.Range("AV2").Value = "cash"
$ws.Range("AW2").Value = 1
$ws.Range("AX2").Value = 0
$ws.Range("BC2").Value = 0
$ws.Range("BG2").Value = 0

# --- Turn the e-mail address into a clickable hyperlink -------------------
$ws.Hyperlinks.Add($ws.Range("G2"), "mailto:ksmith@email.com")
